# 自动更新Excel文件 - 2026-01-10 23:13:10
# 每日巡检：“剩余”天数（E列）递减 1；当“剩余”降为 1（即当天用完）时，
# 重新开始下一周期：剩余 重置为 10，并将“开始时间”（F列，YYYYMMDD 整数）顺延 10 天。
# 若“开始时间”不是规范的 8 位日期（数据异常行），当天跳过、不做任何改动。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string][int]$fVal

    # 跳过开始时间不是标准 8 位 YYYYMMDD 的异常行（例如脏数据 "202510929"）
    if ($fStr.Length -ne 8) {
        continue
    }

    if ([int]$eVal -eq 1) {
        $eCell.Value2 = 10
        $fCell.Value2 = [int]$fVal + 10
    } else {
        $eCell.Value2 = [int]$eVal - 1
    }
}
